$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.1182
$ws.Range("G2").Value = 0.03160595177013854
$ws.Range("H2").Value = 0.03160595177013854
$ws.Range("I2").Value = -0.027603899435608
$ws.Range("J2").Value = -0.027603899435608
$ws.Range("K2").Value = -0.663
$ws.Range("L2").Value = -0.03401744484350949
$ws.Range("U2").Value = 1.881
$ws.Range("V2").Value = 0.07464285714285714
$ws.Range("W2").Value = -0.02188679245283019
$ws.Range("X2").Value = 0.08805005011870837
$ws.Range("Y2").Value = -0.1099368425715386
$ws.Range("Z2").Value = 0.5487822047022385
$ws.Range("AA2").Value = -0.01543540199193816
$ws.Range("AB2").Value = 0.08577291036330373
$ws.Range("AC2").Value = -0.1012083123552419
$ws.Range("AD2").Value = 1.367
$ws.Range("AF2").Value = 1.367
$ws.Range("AG2").Value = -0.514
$ws.Range("AH2").Value = 0.05145481236119998
$ws.Range("AI2").Value = 0.0434423364159278
$ws.Range("AJ2").Value = -0.02082151826946447
$ws.Range("AK2").Value = -0.01737308186304333
$ws.Range("AL2").Value = 0.141
$ws.Range("AM2").Value = 0.106
$ws.Range("AN2").Value = 2.939784946236559
$ws.Range("AO2").Value = -3.815602836879433
$ws.Range("AP2").Value = -1.105376344086021
$ws.Range("AQ2").Value = -5.075471698113208

# Row 3
$ws.Range("D3").Value = -0.08539999999999999
$ws.Range("G3").Value = 0.07007692307692308
$ws.Range("H3").Value = 0.07007692307692308
$ws.Range("I3").Value = -0.008384615384615384
$ws.Range("J3").Value = -0.008384615384615384
$ws.Range("K3").Value = -0.113
$ws.Range("L3").Value = -0.008692307692307692
$ws.Range("U3").Value = 1.53
$ws.Range("V3").Value = 0.1141791044776119
$ws.Range("W3").Value = -0.007106918238993711
$ws.Range("X3").Value = 0.08473951293888962
$ws.Range("Y3").Value = -0.09184643117788333
$ws.Range("Z3").Value = 0.7105766602896967
$ws.Range("AA3").Value = -0.005957911997813611
$ws.Range("AB3").Value = 0.08465985379923031
$ws.Range("AC3").Value = -0.09061776579704392
$ws.Range("AD3").Value = 0.037
$ws.Range("AF3").Value = 0.037
$ws.Range("AG3").Value = -1.493
$ws.Range("AH3").Value = 0.002753590831286745
$ws.Range("AI3").Value = 0.00235114697845841
$ws.Range("AJ3").Value = -0.1253884269757286
$ws.Range("AK3").Value = -0.105089040613782
$ws.Range("AL3").Value = 0.052
$ws.Range("AM3").Value = 0.01699999999999999
$ws.Range("AN3").Value = 0.07326732673267326
$ws.Range("AO3").Value = -2.096153846153846
$ws.Range("AP3").Value = -2.956435643564356
$ws.Range("AQ3").Value = -6.411764705882355

# Row 4
$ws.Range("D4").Value = -0.151
$ws.Range("G4").Value = -0.04545454545454545
$ws.Range("H4").Value = -0.04545454545454545
$ws.Range("I4").Value = -0.06610169491525424
$ws.Range("J4").Value = -0.06610169491525424
$ws.Range("K4").Value = -0.55
$ws.Range("L4").Value = -0.08474576271186442
$ws.Range("U4").Value = 0.351
$ws.Range("V4").Value = 0.0297457627118644
$ws.Range("W4").Value = -0.03666666666666667
$ws.Range("X4").Value = 0.09136058729852713
$ws.Range("Y4").Value = -0.1280272539651938
$ws.Range("Z4").Value = 0.3768873403019745
$ws.Range("AA4").Value = -0.02491289198606272
$ws.Range("AB4").Value = 0.08688596692737713
$ws.Range("AC4").Value = -0.1117988589134399
$ws.Range("AD4").Value = 1.33
$ws.Range("AF4").Value = 1.33
$ws.Range("AG4").Value = 0.9790000000000001
$ws.Range("AH4").Value = 0.1012947448591013
$ws.Range("AI4").Value = 0.08455181182453909
$ws.Range("AJ4").Value = 0.07661006338524141
$ws.Range("AK4").Value = 0.06365823525586839
$ws.Range("AL4").Value = 0.089
$ws.Range("AM4").Value = 0.089
$ws.Range("AN4").Value = -33.25
$ws.Range("AO4").Value = -4.820224719101124
$ws.Range("AP4").Value = -24.475
$ws.Range("AQ4").Value = -4.820224719101124
